$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "42.312.93"
Set-TextCell "E2" "  -3.04%  "
Set-TextCell "D3" "2.221.21"
Set-TextCell "E3" "  -2.15%  "
Set-TextCell "E4" "  +0.27%  "
Set-TextCell "D5" "110.62"
Set-TextCell "E5" "  -8.30%  "
Set-TextCell "D6" "289.84"
Set-TextCell "E6" "  +7.75%  "
Set-TextCell "E7" "  -2.69%  "
Set-TextCell "E8" "  +0.07%  "
Set-TextCell "D9" "0.597"
Set-TextCell "E9" "  -3.73%  "
Set-TextCell "E10" "  -8.27%  "
Set-TextCell "D11" "0.0907"
Set-TextCell "E11" "  -3.81%  "
Set-TextCell "D12" "54.18"
Set-TextCell "E12" "  -0.04%  "
Set-TextCell "E13" "  -8.37%  "
Set-TextCell "D14" "1.00"
Set-TextCell "E14" "  +10.44%  "
Set-TextCell "E15" "  -3.20%  "
Set-TextCell "D16" "14.85"
Set-TextCell "E16" "  -5.77%  "
Set-TextCell "D17" "2.556.20"
Set-TextCell "E17" "  -2.08%  "
Set-TextCell "D18" "2.213.51"
Set-TextCell "E18" "  -2.34%  "
Set-TextCell "D19" "42.290.73"
Set-TextCell "E19" "  -3.00%  "
Set-TextCell "E20" "  +2.90%  "
Set-TextCell "D21" "0.0000104"
Set-TextCell "E21" "  -4.55%  "
Set-TextCell "D22" "72.60"
Set-TextCell "E22" "  -0.08%  "
Set-TextCell "E23" "  +13.18%  "
Set-TextCell "D24" "2.41"
Set-TextCell "E24" "  +0.94%  "
Set-TextCell "D25" "232.04"
Set-TextCell "E25" "  -1.24%  "
Set-TextCell "D26" "8.93"
Set-TextCell "E26" "  -8.28%  "
Set-TextCell "E27" "  -1.58%  "
Set-TextCell "D28" "11.37"
Set-TextCell "E28" "  -7.76%  "
Set-TextCell "E29" "  -2.47%  "
Set-TextCell "D30" "37.62"
Set-TextCell "E30" "  -10.36%  "
Set-TextCell "D31" "173.16"
Set-TextCell "D32" "3.07"
Set-TextCell "E32" "  -8.28%  "
Set-TextCell "D33" "20.83"
Set-TextCell "E33" "  -3.25%  "
Set-TextCell "D34" "0.0877"
Set-TextCell "E34" "  -4.42%  "
Set-TextCell "E35" "  -2.37%  "
Set-TextCell "E36" "  +4.65%  "
Set-TextCell "D37" "4.22"
Set-TextCell "E37" "  -6.11%  "
Set-TextCell "E38" "  -3.64%  "
Set-TextCell "E39" "  -2.45%  "
Set-TextCell "E40" "  -4.86%  "
Set-TextCell "B41" "LidoDAOToken"
Set-TextCell "C41" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D41" "2.38"
Set-TextCell "E41" "  -6.97%  "
Set-TextCell "B42" "MultiversX"
Set-TextCell "C42" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell "D42" "73.02"
Set-TextCell "E42" "  +0.67%  "
Set-TextCell "E43" "  -5.44%  "
Set-TextCell "E44" "  -0.04%  "
Set-TextCell "D45" "12.23"
Set-TextCell "E45" "  -11.31%  "
Set-TextCell "E46" "  -4.52%  "
Set-TextCell "D47" "5.31"
Set-TextCell "E47" "  -7.09%  "
Set-TextCell "E48" "  -0.29%  "
Set-TextCell "D49" "1.65"
Set-TextCell "E49" "  +2.30%  "
Set-TextCell "D50" "101.10"
Set-TextCell "E50" "  -1.62%  "
Set-TextCell "D51" "8.39"
Set-TextCell "E51" "  -2.26%  "
